# Update countries & provincias Spain
#
# This script reproduces the effective content change described by the
# diff: two pairs of country names were reordered in the shared-string
# table (which, combined with the row data, makes the table rows for
# Oman/Croacia, Armenia/Islandia and Letonia/Republica de Chipre swap
# places) and a handful of numeric statistics were refreshed for the
# affected rows (46, 66, 67, 70, 71, 82, 92, 93, 139).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 46 (Australia) : refreshed totals -------------------------------
$ws.Range("B46").Value = 6720
$ws.Range("C46").Value = 4
$ws.Range("D46").Value = 5586
$ws.Range("E46").Value = 1051

# --- Rows 66/67 : Croacia <-> Oman swap places, with refreshed data ------
$ws.Range("A66").Value = "Oman"
$ws.Range("B66").Value = 2049
$ws.Range("C66").Value = 51
$ws.Range("D66").Value = 364
$ws.Range("E66").Value = 1675
$ws.Range("F66").Value = 3
$ws.Range("G66").Value = 0
$ws.Range("H66").Value = 10

$ws.Range("A67").Value = "Croacia"
$ws.Range("B67").Value = 2030
$ws.Range("C67").Value = 0
$ws.Range("D67").Value = 1103
$ws.Range("E67").Value = 872
$ws.Range("F67").Value = 23
$ws.Range("G67").Value = 0
$ws.Range("H67").Value = 55

# --- Rows 70/71 : Islandia <-> Armenia swap places, with refreshed data --
$ws.Range("A70").Value = "Armenia"
$ws.Range("B70").Value = 1808
$ws.Range("C70").Value = 62
$ws.Range("D70").Value = 848
$ws.Range("E70").Value = 931
$ws.Range("F70").Value = 10
$ws.Range("G70").Value = 1
$ws.Range("H70").Value = 29

$ws.Range("A71").Value = "Islandia"
$ws.Range("B71").Value = 1792
$ws.Range("C71").Value = 0
$ws.Range("D71").Value = 1608
$ws.Range("E71").Value = 174
$ws.Range("F71").Value = 3
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 10

# --- Row 82 (Eslovaquia) : refreshed totals -------------------------------
$ws.Range("B82").Value = 1381
$ws.Range("C82").Value = 2
$ws.Range("D82").Value = 403
$ws.Range("E82").Value = 960
$ws.Range("F82").Value = 7

# --- Rows 92/93 : Republica de Chipre <-> Letonia swap, refreshed data ---
$ws.Range("A92").Value = "Letonia"
$ws.Range("B92").Value = 818
$ws.Range("C92").Value = 6
$ws.Range("D92").Value = 267
$ws.Range("E92").Value = 538
$ws.Range("F92").Value = 5
$ws.Range("G92").Value = 1
$ws.Range("H92").Value = 13

$ws.Range("A93").Value = "Republica de Chipre"
$ws.Range("B93").Value = 817
$ws.Range("C93").Value = 0
$ws.Range("D93").Value = 148
$ws.Range("E93").Value = 655
$ws.Range("F93").Value = 15
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 14

# --- Row 139 (Brunei) : refreshed totals ----------------------------------
$ws.Range("D139").Value = 124
$ws.Range("E139").Value = 13

$wb.Save()
